# Updates the "cryptos" price/volume table with freshly scraped figures.
# For price cells whose new text looks like a plain number (e.g. "562.97"),
# the cell is pre-formatted as Text ("@") before the value is written so
# Excel keeps it as a string instead of silently converting it to a number
# (matching the original data, which stores these as text). The number
# format/style is then reset back to the sheet default afterwards so no
# stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.986.00'
$ws.Cells.Item(2, 5).Value = '  -0.32%  '
$ws.Cells.Item(3, 4).Value = '2.419.64'
$ws.Cells.Item(3, 5).Value = '  -0.08%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '562.97'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.09%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '143.14'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.02%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.531'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.26%  '
$ws.Cells.Item(9, 5).Value = '  -0.33%  '
$ws.Cells.Item(11, 5).Value = '  -4.12%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.349'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.11%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.21'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.40%  '
$ws.Cells.Item(14, 5).Value = '  -1.62%  '
$ws.Cells.Item(16, 4).Value = '61.893.21'
$ws.Cells.Item(16, 5).Value = '  -0.08%  '
$ws.Cells.Item(17, 4).Value = '2.417.83'
$ws.Cells.Item(17, 5).Value = '  -0.44%  '
$ws.Cells.Item(18, 5).Value = '  +1.15%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '323.28'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.43%  '
$ws.Cells.Item(20, 5).Value = '  +0.98%  '
$ws.Cells.Item(21, 5).Value = '  -1.21%  '
$ws.Cells.Item(22, 5).Value = '  +0.00%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '66.73'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.98%  '
$ws.Cells.Item(24, 5).Value = '  +0.65%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.79'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.69%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '553.66'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -5.79%  '
$ws.Cells.Item(27, 4).Value = '2.539.24'
$ws.Cells.Item(27, 5).Value = '  -0.12%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.91%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0934'
$ws.Cells.Item(29, 5).Value = '  -0.93%  '
$ws.Cells.Item(30, 5).Value = '  -0.74%  '
$ws.Cells.Item(31, 5).Value = '  -4.52%  '
$ws.Cells.Item(32, 5).Value = '  -1.70%  '
$ws.Cells.Item(33, 5).Value = '  +0.10%  '
$ws.Cells.Item(34, 5).Value = '  -4.00%  '
$ws.Cells.Item(35, 5).Value = '  -0.04%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.73'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.80%  '
$ws.Cells.Item(37, 5).Value = '  -1.56%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '153.09'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.86%  '
$ws.Cells.Item(39, 5).Value = '  -5.00%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.55'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.77%  '
$ws.Cells.Item(41, 5).Value = '  -0.28%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '147.48'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -2.07%  '
$ws.Cells.Item(44, 5).Value = '  -6.04%  '
$ws.Cells.Item(45, 5).Value = '  -0.24%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0528'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -1.98%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.81'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.70%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.592'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.02%  '
$ws.Cells.Item(49, 5).Value = '  -0.52%  '
$ws.Cells.Item(50, 5).Value = '  -0.70%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.56'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.65%  '
